{"js": "// Office.js (Word JavaScript API) edit script.\n// Applies the diff: updates the date line and the 25 multiplication\n// expressions inside the practice-sheet table via exact text\n// search-and-replace, one-to-one (old text -> new text). Every source\n// string below appears exactly once in the document, so a single\n// search + replace per pair reproduces the diff exactly.\nconst body = context.document.body;\n\nconst replacements = [\n  [\"2024-11-30 Saturday\", \"2024-12-01 Sunday\"],\n  [\"267\u00d72=\", \"616\u00d72=\"],\n  [\"212\u00d77=\", \"616\u00d75=\"],\n  [\"774\u00d78=\", \"393\u00d74=\"],\n  [\"177\u00d77=\", \"306\u00d74=\"],\n  [\"763\u00d73=\", \"711\u00d78=\"],\n  [\"229\u00d75=\", \"778\u00d78=\"],\n  [\"712\u00d76=\", \"595\u00d79=\"],\n  [\"408\u00d76=\", \"254\u00d76=\"],\n  [\"586\u00d78=\", \"733\u00d75=\"],\n  [\"455\u00d73=\", \"513\u00d72=\"],\n  [\"241\u00d74=\", \"148\u00d76=\"],\n  [\"788\u00d74=\", \"559\u00d73=\"],\n  [\"205\u00d75=\", \"892\u00d73=\"],\n  [\"619\u00d76=\", \"316\u00d77=\"],\n  [\"308\u00d72=\", \"648\u00d78=\"],\n  [\"748\u00d73=\", \"728\u00d77=\"],\n  [\"933\u00d77=\", \"907\u00d74=\"],\n  [\"150\u00d74=\", \"108\u00d75=\"],\n  [\"824\u00d75=\", \"263\u00d79=\"],\n  [\"370\u00d75=\", \"931\u00d73=\"],\n  [\"300\u00d73=\", \"574\u00d72=\"],\n  [\"265\u00d74=\", \"191\u00d77=\"],\n  [\"223\u00d78=\", \"317\u00d72=\"],\n  [\"533\u00d72=\", \"292\u00d76=\"],\n  [\"485\u00d79=\", \"165\u00d79=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const rng of results.items) {\n    rng.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Applies the diff: updates the date line and the 25 multiplication\n# expressions inside the practice-sheet table via Find/Replace,\n# one-to-one (old text -> new text). Every source string below appears\n# exactly once in the document, so a single Find.Execute replace-all\n# per pair reproduces the diff exactly without touching anything else.\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$pairs = @(\n  @(\"2024-11-30 Saturday\", \"2024-12-01 Sunday\"),\n  @(\"267\u00d72=\", \"616\u00d72=\"),\n  @(\"212\u00d77=\", \"616\u00d75=\"),\n  @(\"774\u00d78=\", \"393\u00d74=\"),\n  @(\"177\u00d77=\", \"306\u00d74=\"),\n  @(\"763\u00d73=\", \"711\u00d78=\"),\n  @(\"229\u00d75=\", \"778\u00d78=\"),\n  @(\"712\u00d76=\", \"595\u00d79=\"),\n  @(\"408\u00d76=\", \"254\u00d76=\"),\n  @(\"586\u00d78=\", \"733\u00d75=\"),\n  @(\"455\u00d73=\", \"513\u00d72=\"),\n  @(\"241\u00d74=\", \"148\u00d76=\"),\n  @(\"788\u00d74=\", \"559\u00d73=\"),\n  @(\"205\u00d75=\", \"892\u00d73=\"),\n  @(\"619\u00d76=\", \"316\u00d77=\"),\n  @(\"308\u00d72=\", \"648\u00d78=\"),\n  @(\"748\u00d73=\", \"728\u00d77=\"),\n  @(\"933\u00d77=\", \"907\u00d74=\"),\n  @(\"150\u00d74=\", \"108\u00d75=\"),\n  @(\"824\u00d75=\", \"263\u00d79=\"),\n  @(\"370\u00d75=\", \"931\u00d73=\"),\n  @(\"300\u00d73=\", \"574\u00d72=\"),\n  @(\"265\u00d74=\", \"191\u00d77=\"),\n  @(\"223\u00d78=\", \"317\u00d72=\"),\n  @(\"533\u00d72=\", \"292\u00d76=\"),\n  @(\"485\u00d79=\", \"165\u00d79=\")\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n\n  $ok = $find.Execute([ref]$oldText, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]$wdFindContinue, [ref]$false, [ref]$newText, [ref]$wdReplaceAll)\n\n  if (-not $ok) {\n    Write-Output (\"Replacement failed for: \" + $oldText)\n  }\n}\n\nWrite-Output \"done\"\n"}
